$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value into a cell while forcing text storage
# (avoids Excel auto-converting numeric-looking strings like "305.99"
# into numbers), and resets the style afterward so no stray number-format
# style index is left attached to the cell.
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "42.982.73"
Set-TextCell $ws.Range("E2") "  +0.09%  "

# Row 3
Set-TextCell $ws.Range("D3") "2.537.45"
Set-TextCell $ws.Range("E3") "  -1.01%  "

# Row 4
Set-TextCell $ws.Range("E4") "  +0.03%  "

# Row 5
Set-TextCell $ws.Range("D5") "305.99"
Set-TextCell $ws.Range("E5") "  +1.29%  "

# Row 6
Set-TextCell $ws.Range("D6") "101.18"
Set-TextCell $ws.Range("E6") "  +7.13%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.580"
Set-TextCell $ws.Range("E7") "  +1.13%  "

# Row 8
Set-TextCell $ws.Range("E8") "  +0.03%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.549"
Set-TextCell $ws.Range("E9") "  +0.72%  "

# Row 10
Set-TextCell $ws.Range("D10") "37.60"
Set-TextCell $ws.Range("E10") "  +3.99%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.0820"
Set-TextCell $ws.Range("E11") "  +1.21%  "

# Row 12
Set-TextCell $ws.Range("D12") "7.63"
Set-TextCell $ws.Range("E12") "  -1.68%  "

# Row 14
Set-TextCell $ws.Range("D14") "2.925.61"
Set-TextCell $ws.Range("E14") "  -0.98%  "

# Row 15
Set-TextCell $ws.Range("B15") "WrappedEther"
Set-TextCell $ws.Range("C15") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws.Range("D15") "2.541.36"
Set-TextCell $ws.Range("E15") "  -0.56%  "

# Row 16
Set-TextCell $ws.Range("B16") "Chainlink"
Set-TextCell $ws.Range("C16") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws.Range("D16") "15.27"
Set-TextCell $ws.Range("E16") "  +7.41%  "

# Row 17
Set-TextCell $ws.Range("E17") "  -0.60%  "

# Row 18
Set-TextCell $ws.Range("D18") "42.979.35"
Set-TextCell $ws.Range("E18") "  -0.03%  "

# Row 19
Set-TextCell $ws.Range("D19") "13.18"
Set-TextCell $ws.Range("E19") "  +3.24%  "

# Row 20
Set-TextCell $ws.Range("D20") "0.0₃0989"
Set-TextCell $ws.Range("E20") "  -0.65%  "

# Row 21
Set-TextCell $ws.Range("D21") "6.51"
Set-TextCell $ws.Range("E21") "  -0.69%  "

# Row 22
Set-TextCell $ws.Range("D22") "71.73"
Set-TextCell $ws.Range("E22") "  +0.35%  "

# Row 23
Set-TextCell $ws.Range("D23") "254.43"
Set-TextCell $ws.Range("E23") "  +0.49%  "

# Row 24
Set-TextCell $ws.Range("E24") "  -0.29%  "

# Row 25
Set-TextCell $ws.Range("E25") "  -3.06%  "

# Row 26
Set-TextCell $ws.Range("D26") "27.33"
Set-TextCell $ws.Range("E26") "  -4.79%  "

# Row 27
Set-TextCell $ws.Range("E27") "  +0.20%  "

# Row 28
Set-TextCell $ws.Range("E28") "  +9.38%  "

# Row 29
Set-TextCell $ws.Range("D29") "10.43"
Set-TextCell $ws.Range("E29") "  +1.75%  "

# Row 30
Set-TextCell $ws.Range("D30") "38.77"
Set-TextCell $ws.Range("E30") "  +4.42%  "

# Row 31
Set-TextCell $ws.Range("E31") "  +1.49%  "

# Row 32
Set-TextCell $ws.Range("D32") "158.13"
Set-TextCell $ws.Range("E32") "  +3.05%  "

# Row 33
Set-TextCell $ws.Range("E33") "  -1.38%  "

# Row 34
Set-TextCell $ws.Range("D34") "0.0800"
Set-TextCell $ws.Range("E34") "  +0.09%  "

# Row 35
Set-TextCell $ws.Range("D35") "3.30"
Set-TextCell $ws.Range("E35") "  -2.15%  "

# Row 36
Set-TextCell $ws.Range("E36") "  -3.72%  "

# Row 37
Set-TextCell $ws.Range("D37") "18.48"
Set-TextCell $ws.Range("E37") "  +3.29%  "

# Row 38
Set-TextCell $ws.Range("D38") "0.116"
Set-TextCell $ws.Range("E38") "  +1.57%  "

# Row 39
Set-TextCell $ws.Range("D39") "0.120"
Set-TextCell $ws.Range("E39") "  +0.33%  "

# Row 40
Set-TextCell $ws.Range("D40") "23.82"
Set-TextCell $ws.Range("E40") "  +3.17%  "

# Row 41
Set-TextCell $ws.Range("E41") "  +2.23%  "

# Row 42
Set-TextCell $ws.Range("D42") "2.09"
Set-TextCell $ws.Range("E42") "  +3.01%  "

# Row 43
Set-TextCell $ws.Range("D43") "3.87"
Set-TextCell $ws.Range("E43") "  +0.06%  "

# Row 44
Set-TextCell $ws.Range("E44") "  -1.71%  "

# Row 45
Set-TextCell $ws.Range("D45") "0.998"
Set-TextCell $ws.Range("E45") "  -0.05%  "

# Row 46
Set-TextCell $ws.Range("D46") "2.045.89"
Set-TextCell $ws.Range("E46") "  -2.60%  "

# Row 47
Set-TextCell $ws.Range("D47") "86.27"
Set-TextCell $ws.Range("E47") "  +1.45%  "

# Row 48
Set-TextCell $ws.Range("D48") "9.01"
Set-TextCell $ws.Range("E48") "  -2.18%  "

# Row 49
Set-TextCell $ws.Range("D49") "2.781.87"
Set-TextCell $ws.Range("E49") "  -0.92%  "

# Row 50
Set-TextCell $ws.Range("D50") "0.194"
Set-TextCell $ws.Range("E50") "  +1.50%  "

# Row 51
Set-TextCell $ws.Range("D51") "103.56"
Set-TextCell $ws.Range("E51") "  -2.77%  "
